$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.109331011772156
$ws.Range("B1").Value = 1.946161389350891
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.035928726196289
$ws.Range("E1").Value = 1.112438440322876
